$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules Sheet")

# Insert a new row before row 14; this pushes the existing row 14
# ("Please Do See Sample Sheet give below") down to row 15, preserving
# its A14:D14 merge as A15:D15.
$ws.Rows.Item(14).Insert()

# Row 13 gains a sequence number and the new validation message.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Only One Sheet is Allowed. Otherwise Data cannot be Saved."

# Make the Rules Sheet the active tab (matches workbook.xml activeTab change).
$ws.Activate()

$wb.Save()
